# "Generate Report for Handoff" - a new handoff was generated for b.md,
# moving its status from "Handed back: in sync with en-US" to
# "Ready for handoff" and recording the new handoff artifacts/timestamps.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e9251d2fd2344319e672556d3cbdbbad0eeeea6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4b0e343644d5c4dd6d02fb2aec7ead6b89a90149/e2e/b.md."

# --- Overview sheet: row 3 corresponds to b.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-27 02:37:06"

# --- zh-cn sheet: row 3 corresponds to b.md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$savedStyle = $wsZhCn.Range("F3").Style
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = $savedStyle
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-27 02:36:58"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 corresponds to b.md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$savedStyle2 = $wsDeDe.Range("F3").Style
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = $savedStyle2
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-27 02:37:06"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
